# Generate Report for Handoff
# Rotate the localization-status report for a new run: new GUID-named
# source file, new xliff content hashes, and refreshed handoff/handback
# timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "615f9ee1-53bb-403d-8c17-b4c0331583d6"
$newGuid = "88397580-d60d-44c7-8f0a-16fb571530b6"

$oldHash = "0ebe14d94cc26a57afdb471f3cbf62791ec2fe28"
$newHash = "62cded27da37696ec18b99c658c58f586eccb7f4"

$newGenerateDate  = "2016-08-21 15:04:33"
$newHandoffDate   = "2016-08-21 15:04:28"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newPathDisplay = "e2e\$newGuid.md"
$oldLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32c0cadbec4b96ee5da40c3dabe2e26de5fa0051/e2e/$oldGuid.md"

# Pre-existing (file-loaded) hyperlinks can't be edited/removed in place
# through Hyperlinks.Item(n) -- only the Range-level Hyperlinks collection
# actually detaches them. Drop it and re-add with the same target address
# (r:id / URL is unchanged per the diff) but the refreshed display text.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldLinkAddress, "", "", $newPathDisplay)

$wsOverview.Range("G2").Value = $newGenerateDate

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"

$oldSourceLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32c0cadbec4b96ee5da40c3dabe2e26de5fa0051/e2e/$oldGuid.md"
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $oldSourceLinkAddress, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = $newHandoffDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $oldSourceLinkAddress, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newGenerateDate
